$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the number formats currently used:
#  - A2:A51 use the regular date/time format (same as A51)
#  - A52 (old last row) uses the special date-only "last row" format
$normalDateFormat = $ws.Range("A51").NumberFormat
$lastRowDateFormat = $ws.Range("A52").NumberFormat

# Row 52 is no longer the last row, so give it the regular date/time format
$ws.Range("A52").NumberFormat = $normalDateFormat

# Add new row 53 with data
$ws.Range("A53").Value = 45638
$ws.Range("B53").Value = 130
$ws.Range("C53").Value = 116
$ws.Range("D53").Value = 120

# Row 53 is now the last row: apply the special "last row" date format
$ws.Range("A53").NumberFormat = $lastRowDateFormat
